# Apply "search" roll number value to all data rows in column B,
# and move the active cell selection to B4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update RollNumber (column B) for all data rows (2-4) to the new
# shared value, simulating a search parameter passed to the controller.
$ws.Range("B2").Value = "HE130576"
$ws.Range("B3").Value = "HE130576"
$ws.Range("B4").Value = "HE130576"

# Update the active selection to B4.
$ws.Range("B4").Select()
